$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "先导智能"
$ws.Range("A3").Value = "航天电子"
$ws.Range("C3").Value = "锋龙股份"
$ws.Range("A4").Value = "隆基绿能"
$ws.Range("B4").Value = "利欧股份"
$ws.Range("A5").Value = "钧达股份"
$ws.Range("B5").Value = "拓日新能"
$ws.Range("C5").Value = "嘉美包装"
$ws.Range("A6").Value = "乾照光电"
$ws.Range("B6").Value = "协鑫集成"
$ws.Range("C6").Value = "航天电子"
$ws.Range("A7").Value = "拓日新能"
$ws.Range("A8").Value = "利欧股份"
$ws.Range("B8").Value = "航天电子"
$ws.Range("C8").Value = "隆基绿能"
$ws.Range("A9").Value = "中国卫星"
$ws.Range("B9").Value = "钧达股份"
$ws.Range("C9").Value = "白银有色"
$ws.Range("A10").Value = "协鑫集成"
$ws.Range("B10").Value = "白银有色"
$ws.Range("C10").Value = "湖南白银"
$ws.Range("A11").Value = "锋龙股份"
$ws.Range("C11").Value = "协鑫集成"
$ws.Range("A12").Value = "白银有色"
$ws.Range("B12").Value = "中国卫星"
$ws.Range("C12").Value = "东方日升"
$ws.Range("A13").Value = "金风科技"
$ws.Range("B13").Value = "锋龙股份"
$ws.Range("C13").Value = "拓日新能"
$ws.Range("A14").Value = "迈为股份"
$ws.Range("B14").Value = "铜陵有色"
$ws.Range("C14").Value = "钧达股份"
$ws.Range("A15").Value = "嘉美包装"
$ws.Range("B15").Value = "捷佳伟创"
$ws.Range("C15").Value = "乾照光电"
$ws.Range("A16").Value = "先导智能"
$ws.Range("B16").Value = "先导智能"
$ws.Range("C16").Value = "蓝色光标"
$ws.Range("A17").Value = "湖南白银"
$ws.Range("B17").Value = "湖南白银"
$ws.Range("C17").Value = "岩山科技"
$ws.Range("A18").Value = "巨力索具"
$ws.Range("B18").Value = "嘉美包装"
$ws.Range("C18").Value = "中国卫星"
$ws.Range("A19").Value = "浙文互联"
$ws.Range("B19").Value = "迈为股份"
$ws.Range("C19").Value = "巨力索具"
$ws.Range("B20").Value = "贵州茅台"
$ws.Range("C20").Value = "通富微电"
$ws.Range("A21").Value = "蓝色光标"
$ws.Range("B21").Value = "东方财富"
$ws.Range("C21").Value = "中国长城"
